$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: "<Basename>_old" -> "<Basename>_FV2304" and
#    "<Basename>_new" -> "<Basename>_FV2310" (the "diff" header in K1 stays).
$headers = @{
  "A1" = "Segmentname_FV2304";        "B1" = "Segmentgruppe_FV2304";
  "C1" = "Segment_FV2304";            "D1" = "Datenelement_FV2304";
  "E1" = "Segment ID_FV2304";         "F1" = "Code_FV2304";
  "G1" = "Qualifier_FV2304";          "H1" = "Beschreibung_FV2304";
  "I1" = "Bedingungsausdruck_FV2304"; "J1" = "Bedingung_FV2304";
  "L1" = "Segmentname_FV2310";        "M1" = "Segmentgruppe_FV2310";
  "N1" = "Segment_FV2310";            "O1" = "Datenelement_FV2310";
  "P1" = "Segment ID_FV2310";         "Q1" = "Code_FV2310";
  "R1" = "Qualifier_FV2310";          "S1" = "Beschreibung_FV2310";
  "T1" = "Bedingungsausdruck_FV2310"; "U1" = "Bedingung_FV2310"
}
foreach ($addr in $headers.Keys) {
  $ws.Range($addr).Value = $headers[$addr]
}

# 2) Turn the used range into an Excel Table ("Table1") so the header row
#    drives an auto filter, matching the new sharedStrings-backed headers.
$range = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

# 3) Freeze the header row (split after row 1, active pane bottom-left).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
